# Add two new columns (I: "I0", J: "IF") to Sheet1, mirroring the
# existing header style used by the other header cells (B1:H1).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Headers in row 1
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

# Match the style of the existing header cells (e.g. H1) for the new headers
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122) # xlPasteFormats

# Data values for columns I and J, rows 2-15
$data = @{
    2  = @(1, 3)
    3  = @(1, 3)
    4  = @(8, 8)
    5  = @(8, 9)
    6  = @(4, 5)
    7  = @(7, 8)
    8  = @(7, 8)
    9  = @(6, 8)
    10 = @(6, 7)
    11 = @(3, 6)
    12 = @(2, 4)
    13 = @(1, 5)
    14 = @(6, 7)
    15 = @(7, 8)
}

foreach ($row in $data.Keys) {
    $vals = $data[$row]
    $ws.Cells.Item($row, 9).Value = $vals[0]   # column I
    $ws.Cells.Item($row, 10).Value = $vals[1]  # column J
}
